$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.467.00"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.880.78"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'240.58"
$ws.Range("E5").Value = "  +3.78%  "
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "'42.92"
$ws.Range("E8").Value = "  +7.28%  "
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("D10").Value = "'0.0700"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").Value = "'0.0990"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "2.151.24"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").Value = "1.927.26"
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("D14").Value = "'11.67"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "35.457.79"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "'70.76"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").Value = "'242.54"
$ws.Range("D21").Value = "'12.40"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'2.27"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'170.47"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "'1.91"
$ws.Range("E26").Value = "  +25.50%  "
$ws.Range("E27").Value = "  +5.90%  "
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'1.82"
$ws.Range("E33").Value = "  +24.38%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.06"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("D35").Value = "'2.05"
$ws.Range("E35").Value = "  +7.32%  "
$ws.Range("D36").Value = "'0.826"
$ws.Range("E36").Value = "  +18.44%  "
$ws.Range("E37").Value = "  +7.23%  "
$ws.Range("D38").Value = "'1.11"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("E39").Value = "  +5.30%  "
$ws.Range("D40").Value = "'91.37"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").Value = "1.355.89"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").Value = "'15.26"
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("E43").Value = "  +15.44%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.37"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("B45").Value = "Gas"
$ws.Range("C45").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D45").Value = "'13.06"
$ws.Range("E45").Value = "  +56.76%  "
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E47").Value = "  +6.57%  "
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("D49").Value = "2.064.42"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("D50").Value = "'0.0690"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").Value = "'3.45"
$ws.Range("E51").Value = "  -0.18%  "
